$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.18%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.65%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.20%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08045"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.14%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.879"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.26%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.795"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.59%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9257"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.32%"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1404"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.98%"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1901"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.23%"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09088"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.12%"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03435"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.80%"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09889"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.07%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001400"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.11%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006039"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-9.24%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.843"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.49%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.123"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.94%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "13.41%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3421"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.10%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1322"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.18%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.807"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.15%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2613"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.97%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04350"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.57%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.52%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004293"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.85%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001297"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.19%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "42.16%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02011"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.57%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05116"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007514"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.46%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01008"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.32%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1356"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.68%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002156"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.67%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009624"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.40%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006236"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.16%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.20%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.73"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.26%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001248"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-22.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.20%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.20%"
